# Feat: Atualização da Classificação de Defeitos
#
# Adds a new "RESPONSABILIDADE" column (C) to the catalog, renames the
# header row, trims the leading space from the descriptive labels in
# column B, adds a duplicate "ENG" row, and appends a trailer row
# marked "NÃO MOSTRAR NO ÍNDICE".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A, C, D new/changed; B unchanged) ---
$ws.Columns.Item(1).ColumnWidth = 23.69
$ws.Columns.Item(2).ColumnWidth = 34.25
$ws.Columns.Item(3).ColumnWidth = 24.11
$ws.Columns.Item(4).ColumnWidth = 20.22

# --- Header row ---
$ws.Range("A1").Value = "CÓDIGO DO FORNECEDOR"
$ws.Range("B1").Value = "CLASSIFICAÇÃO DO FORNECEDOR"
$ws.Range("C1").Value = "RESPONSABILIDADE"

# --- Data rows ---
$ws.Range("A2").Value = "F"
$ws.Range("B2").Value = "CHINA"
$ws.Range("C2").Value = "CHINA"

$ws.Range("A3").Value = "FL "
$ws.Range("B3").Value = "FORNECEDOR LOCAL"
$ws.Range("C3").Value = "FORNECEDOR LOCAL"

$ws.Range("A4").Value = "JIG "
$ws.Range("B4").Value = "ENGENHARIA/PROJETO"
$ws.Range("C4").Value = "ENGENHARIA/PROJETO"

$ws.Range("A5").Value = "ENG "
$ws.Range("B5").Value = "ENGENHARIA/PROJETO"
$ws.Range("C5").Value = "ENGENHARIA/PROJETO"

$ws.Range("A6").Value = "IP "
$ws.Range("B6").Value = "PROCESSO INJEÇÃO"
$ws.Range("C6").Value = "PROCESSO INJEÇÃO"

$ws.Range("A7").Value = "LCM "
$ws.Range("B7").Value = "PROCESSO LCM"
$ws.Range("C7").Value = "PROCESSO LCM"

$ws.Range("A8").Value = "MA "
$ws.Range("B8").Value = "PROCESSO MA"
$ws.Range("C8").Value = "PROCESSO MA"

$ws.Range("A9").Value = "AF "
$ws.Range("B9").Value = "PROC. ALTO FALANTE"
$ws.Range("C9").Value = "PROC. ALTO FALANTE"

$ws.Range("A10").Value = "DP "
$ws.Range("B10").Value = "DIP PTH"
$ws.Range("C10").Value = "DIP PTH"

$ws.Range("A11").Value = "T "
$ws.Range("B11").Value = "PROCESSO PTH"
$ws.Range("C11").Value = "PROCESSO PTH"

$ws.Range("A12").Value = "P "
$ws.Range("B12").Value = "PROCESSO PA"
$ws.Range("C12").Value = "PROCESSO PA"

$ws.Range("A13").Value = "AC "
$ws.Range("B13").Value = "ACÚMULO"
$ws.Range("C13").Value = "NÃO MOSTRAR NO ÍNDICE"

$ws.Range("A14").Value = "OC "
$ws.Range("B14").Value = "OCORRÊNCIA"
$ws.Range("C14").Value = "NÃO MOSTRAR NO ÍNDICE"

$ws.Range("A15").Value = "RT "
$ws.Range("B15").Value = "RETRABALHO PEÇA"
$ws.Range("C15").Value = "NÃO MOSTRAR NO ÍNDICE"

$ws.Range("A16").Value = "INT MOD "
$ws.Range("B16").Value = "INTRODUÇÃO DE MODELO"
$ws.Range("C16").Value = "NÃO MOSTRAR NO ÍNDICE"

$ws.Range("A17").Value = "VER "
$ws.Range("B17").Value = "REVISÃO"
$ws.Range("C17").Value = "NÃO MOSTRAR NO ÍNDICE"

$ws.Range("A18").Value = "FF "
$ws.Range("B18").Value = "FALSA FALHA "
$ws.Range("C18").Value = "NÃO MOSTRAR NO ÍNDICE"

$ws.Range("A19").Value = "AF RET "
$ws.Range("B19").Value = "ALTO FALANTE RETRABALHO"
$ws.Range("C19").Value = "NÃO MOSTRAR NO ÍNDICE"

$ws.Range("A20").Value = "RC "
$ws.Range("B20").Value = "RETORNO DE CAMPO"
$ws.Range("C20").Value = "NÃO MOSTRAR NO ÍNDICE"

$ws.Range("A21").Value = "A "
$ws.Range("B21").Value = "ARMAZENAMENTO"
$ws.Range("C21").Value = "NÃO MOSTRAR NO ÍNDICE"

# --- New trailer row, all three cells marked as not shown in the index ---
$ws.Range("A22").Value = "NÃO MOSTRAR NO ÍNDICE"
$ws.Range("B22").Value = "NÃO MOSTRAR NO ÍNDICE"
$ws.Range("C22").Value = "NÃO MOSTRAR NO ÍNDICE"

# --- Selection moved to D9, matching the authored workbook view ---
$ws.Range("D9").Select()
